# Fix for line 476: num_of_travellers / travellers_change were computed with
# the wrong aggregation (looked like an average/count instead of a SUM over
# the underlying traveller records). This rewrites Sheet1's data body
# (A2:H44) with the corrected "sum" aggregation, which also pulls in rows
# that were previously missing from the (wrongly) aggregated pivot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(106,'February','Germany','EN',189.81,19,1800,1800),
    @(115,'February','Italy','ES',79.92,8,300,300),
    @(221,'March','Austria','EN',49.95,5,400,400),
    @(222,'March','China','EN',39.96,4,1214.473684210527,300),
    @(224,'March','France','EN',388.96,16,194.6889915902719,220),
    @(225,'March','Germany','EN',2037.96,204,973.6842105263157,973.6842105263157),
    @(227,'March','Greece','EN',61245.7,3185,242.6118195671096,205.9558117195005),
    @(232,'March','Italy','DE',2652.88,113,397.2782484816675,351.9999999999999),
    @(246,'March','Spain','DE',1356,65,236.8944099378882,261.1111111111111),
    @(2,'April','Egypt','EN',185.76,15,382.4935064935065,400),
    @(8,'April','Greece','FR',45373.76,1782,633.6939244440348,384.2391304347826),
    @(13,'April','Italy','ES',203.8,20,240.0066733400067,233.3333333333333),
    @(15,'April','Italy','IT',2213.47,92,415.5277622507919,228.5714285714286),
    @(17,'April','Netherlands','EN',73.83,9,824.0300375469337,800),
    @(18,'April','Portugal','DE',4174.2,205,1305.454545454545,1181.25),
    @(19,'April','Portugal','EN',78870.77,4046,219.734137194483,213.4004647560031),
    @(20,'April','Portugal','ES',26678.2,1399,672.5757641573636,632.4607329842931),
    @(21,'April','Portugal','FR',14614.88,789,332.9590766623811,333.5164835164835),
    @(22,'April','Portugal','IT',1193.9,54,421.353711790393,390.9090909090909),
    @(28,'April','Turkey','EN',2209.17,92,2528.399762046401,1214.285714285714),
    @(205,'June','Netherlands','EN',39.96,4,300,300),
    @(215,'June','Spain','ES',1083.77,60,204.0454481694487,233.3333333333333),
    @(187,'July','United States','EN',259,35,1380,1650),
    @(32,'August','China','EN',30.96,4,209.9099099099099,300),
    @(46,'August','Italy','FR',554,25,1267.901234567901,1150),
    @(55,'August','Portugal','IT',4283.9,196,241.8914604948124,216.1290322580645),
    @(62,'August','Spain','FR',3077.81,136,318.3739771089906,300),
    @(64,'August','Turkey','DE',608.5,27,305.3964023984011,285.7142857142857),
    @(65,'August','Turkey','EN',5100.38,148,266.1514605483248,208.3333333333333),
    @(362,'September','Greece','CH',432.5,15,3507.172643869892,1400),
    @(375,'September','Netherlands','DE',39.95,5,400,400),
    @(384,'September','Spain','DE',1978.1,138,247.8834350433513,305.8823529411765),
    @(388,'September','Spain','IT',585.02,52,132.5013909864081,300),
    @(391,'September','United States','EN',419.37,57,324.506529000911,338.4615384615385),
    @(325,'October','Germany','EN',29.97,3,200,200),
    @(326,'October','Greece','CH',1763.47,68,307.7387283236994,353.3333333333333),
    @(331,'October','Greece','GR',322.05,34,386.6273798730734,325),
    @(334,'October','Italy','DE',3461.89,175,344.5159219311762,306.9767441860465),
    @(354,'October','United Kingdom','EN',65.78999999999999,3,558.5585585585585,200),
    @(289,'November','France','EN',2025.98,203,1222.354937667254,2155.555555555556),
    @(290,'November','Germany','EN',10027.37,1004,33358.02469135803,33366.66666666667),
    @(318,'November','United States','EN',2237.02,226,604.1518461393182,629.0322580645161),
    @(86,'December','Netherlands','EN',35.97,3,350.1877346683353,200),
)

$startRow = 2
$oldLastRow = 37
$newLastRow = $startRow + $data.Length - 1

# Copy the column-A style (bold/border/center) from an existing styled cell
# so any newly-added rows below the old data range pick up the same formatting.
$ws.Cells.Item(2, 1).Copy()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    if ($r -gt $oldLastRow) {
        $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    }
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

$excel.CutCopyMode = $false

Write-Host "Done. Rows written:" $data.Length